$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column G first (shifts H:L left to G:K).
# Original row1 G1:L1 = 4,83,34,106,3,1 -> after this delete: G:K = 83,34,106,3,1
$ws.Columns("G").Delete()

# Now the old J1 value (106) sits in column I. Delete that column too
# (shifts J:K left to I:J) -> G:J = 83,34,3,1
$ws.Columns("I").Delete()
